{"js": "// Apply Korean translations to the document body text, and to the\n// single comment, matching the target diff.\n//\n// Strategy: use Range.search() with matchCase to uniquely locate each\n// English source string (including punctuation) and replace it in place\n// with insertText(..., Word.InsertLocation.replace). This preserves the\n// existing run/paragraph formatting of the located text.\n\nasync function replaceOnce(body, searchText, replacement, options) {\n  const searchOptions = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(searchText, searchOptions);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  // Only replace the first match found; callers pass strings that are\n  // unique within the body unless they intentionally want all matches.\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nasync function replaceAll(body, searchText, replacement, options) {\n  const searchOptions = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(searchText, searchOptions);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1 & 3) \"English\" -> \"\uc601\uc5b4\" (appears twice: the hyperlink text, and the\n// plain-run heading further down). Both occurrences are translated.\nawait replaceAll(body, \"English\", \"\uc601\uc5b4\");\n\n// 2) Language list line.\nawait replaceOnce(\n  body,\n  \" / Portuguese / French / Thai / Vietnamese / Spanish\",\n  \" / \ud3ec\ub974\ud22c\uac08\uc5b4 / \ud504\ub791\uc2a4\uc5b4 / \ud0dc\uad6d\uc5b4 / \ubca0\ud2b8\ub0a8\uc5b4 / \uc2a4\ud398\uc778\uc5b4\"\n);\n\n// 4) \"Brief\" label.\nawait replaceOnce(body, \"Brief\", \"\uac1c\uc694\");\n\n// 5) Brief description paragraph (keep the trailing English sentence\n// about customer.io untouched, matching the diff).\nawait replaceOnce(\n  body,\n  \"An email sent upon verification to partners in the target country who have sent the correct documents. It will be sent via customer.io\",\n  \"\uc815\ud655\ud55c \uc11c\ub958\ub97c \uc81c\ucd9c\ud55c \ud604\uc9c0 \ud30c\ud2b8\ub108\uc5d0\uac8c \uac80\ud1a0 \uc644\ub8cc \ud6c4 \ubc1c\uc1a1\ub418\ub294 \uc774\uba54\uc77c\uc785\ub2c8\ub2e4. It will be sent via customer.io\"\n);\n\n// 6) \"Target audience\" label.\nawait replaceOnce(body, \"Target audience\", \"\ub300\uc0c1 \ub3c5\uc790\");\n\n// 7) Target audience description.\nawait replaceOnce(\n  body,\n  \"Invited partners who didn\\u2019t submit their documents on time\",\n  \"\uc81c\ub54c \uc11c\ub958\ub97c \uc81c\ucd9c\ud558\uc9c0 \uc54a\uc740 \ucd08\uccad\ub41c \ud30c\ud2b8\ub108\"\n);\n\n// 8) \"Subject line\" label.\nawait replaceOnce(body, \"Subject line\", \"\uc81c\ubaa9\");\n\n// 9) Subject line tail.\nawait replaceOnce(body, \" \u2014 one step closer!\", \" \u2014 \ud55c \uac78\uc74c \ub354 \uac00\uae4c\uc6cc\uc84c\uc2b5\ub2c8\ub2e4!\");\n\n// 10) Big centered heading.\nawait replaceOnce(\n  body,\n  \"Your documents have been verified!\",\n  \"\uadc0\ud558\uc758 \ubb38\uc11c\uac00 \ud655\uc778\ub418\uc5c8\uc2b5\ub2c8\ub2e4!\"\n);\n\n// 11) Greeting.\nawait replaceOnce(body, \"Hi \", \"\uc548\ub155\ud558\uc138\uc694 \");\n\n// 12) Reminder paragraph.\nawait replaceOnce(\n  body,\n  \"We\\u2019ll be sending out more details about the event soon, including the agenda and travel arrangements, so make sure to check your inbox regularly.\",\n  \"\uc548\uac74\uacfc \uc5ec\ud589 \uc77c\uc815\uc744 \ud3ec\ud568\ud558\uc5ec \ud589\uc0ac\uc5d0 \ub300\ud55c \uc790\uc138\ud55c \ub0b4\uc6a9\uc744 \uace7 \ubcf4\ub0b4\ub4dc\ub9b4 \uc608\uc815\uc774\ub2c8 \ubc1b\uc740 \ud3b8\uc9c0\ud568\uc744 \uc815\uae30\uc801\uc73c\ub85c \ud655\uc778\ud574 \uc8fc\uc2dc\uae30 \ubc14\ub78d\ub2c8\ub2e4.\"\n);\n\n// 13) Country manager contact intro.\nawait replaceOnce(\n  body,\n  \"If you have any questions, please contact your country manager, \",\n  \"\uad81\uae08\ud558\uc2e0 \uc0ac\ud56d\uc740, \uadc0\ud558\uc758 \uad6d\uac00 \ub2f4\ub2f9\uc790 \"\n);\n\n// 14) \", at \" -> \"\uc5d0\uac8c \" (between [NAME] and [EMAIL ADDRESS]).\nawait replaceOnce(body, \", at \", \"\uc5d0\uac8c \");\n\n// 15) \" or \" between [EMAIL ADDRESS] and [WHATSAPP NO]. The earlier\n// \"live chat or WhatsApp\" occurrence of \" or \" must remain untouched, so\n// search only within the paragraph that contains \"[WHATSAPP NO]\".\n{\n  const whatsappNoResults = body.search(\"[WHATSAPP NO]\", { matchCase: true });\n  whatsappNoResults.load(\"items\");\n  await context.sync();\n  const targetPara = whatsappNoResults.items[0].paragraphs.getFirst();\n  const orResults = targetPara.search(\" or \", { matchCase: true });\n  orResults.load(\"items\");\n  await context.sync();\n  orResults.items[0].insertText(\" \ub610\ub294 \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 16) WhatsApp closing sentence.\nawait replaceOnce(body, \" (WhatsApp). \", \" (WhatsApp)\uc744 \ud1b5\ud574 \uc5f0\ub77d\ud574 \uc8fc\uc2dc\uae30 \ubc14\ub78d\ub2c8\ub2e4. \");\n\n// 17) The single comment's text.\nconst comments = context.document.body.getComments();\ncomments.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < comments.items.length; i++) {\n  comments.items[i].load(\"content\");\n}\nawait context.sync();\nfor (let i = 0; i < comments.items.length; i++) {\n  if (comments.items[i].content === \"choose either one\") {\n    comments.items[i].content = \"\ub458 \uc911 \ud558\ub098\ub97c \uc120\ud0dd\ud558\uc138\uc694\";\n  }\n}\nawait context.sync();\n", "ps1": "# Apply Korean translations to the document body text, and to the\n# single comment, matching the target diff.\n#\n# Strategy: use Range.Find.Execute(..., Replace:=wdReplaceAll) against\n# $d.Content (or a scoped sub-range, when a phrase is not unique) so\n# that the existing run/paragraph formatting of the located text is\n# preserved.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($range, [string]$findText, [string]$replaceText) {\n    $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1 & 3) \"English\" -> \"\uc601\uc5b4\" (appears twice: the hyperlink text, and the\n# plain-run heading further down). Both occurrences are translated.\nReplace-All $d.Content \"English\" \"\uc601\uc5b4\"\n\n# 2) Language list line.\nReplace-All $d.Content \" / Portuguese / French / Thai / Vietnamese / Spanish\" \" / \ud3ec\ub974\ud22c\uac08\uc5b4 / \ud504\ub791\uc2a4\uc5b4 / \ud0dc\uad6d\uc5b4 / \ubca0\ud2b8\ub0a8\uc5b4 / \uc2a4\ud398\uc778\uc5b4\"\n\n# 4) \"Brief\" label.\nReplace-All $d.Content \"Brief\" \"\uac1c\uc694\"\n\n# 5) Brief description paragraph (keep the trailing English sentence\n# about customer.io untouched, matching the diff).\nReplace-All $d.Content \"An email sent upon verification to partners in the target country who have sent the correct documents. It will be sent via customer.io\" \"\uc815\ud655\ud55c \uc11c\ub958\ub97c \uc81c\ucd9c\ud55c \ud604\uc9c0 \ud30c\ud2b8\ub108\uc5d0\uac8c \uac80\ud1a0 \uc644\ub8cc \ud6c4 \ubc1c\uc1a1\ub418\ub294 \uc774\uba54\uc77c\uc785\ub2c8\ub2e4. It will be sent via customer.io\"\n\n# 6) \"Target audience\" label.\nReplace-All $d.Content \"Target audience\" \"\ub300\uc0c1 \ub3c5\uc790\"\n\n# 7) Target audience description.\nReplace-All $d.Content \"Invited partners who didn\u2019t submit their documents on time\" \"\uc81c\ub54c \uc11c\ub958\ub97c \uc81c\ucd9c\ud558\uc9c0 \uc54a\uc740 \ucd08\uccad\ub41c \ud30c\ud2b8\ub108\"\n\n# 8) \"Subject line\" label.\nReplace-All $d.Content \"Subject line\" \"\uc81c\ubaa9\"\n\n# 9) Subject line tail.\nReplace-All $d.Content \" \u2014 one step closer!\" \" \u2014 \ud55c \uac78\uc74c \ub354 \uac00\uae4c\uc6cc\uc84c\uc2b5\ub2c8\ub2e4!\"\n\n# 10) Big centered heading.\nReplace-All $d.Content \"Your documents have been verified!\" \"\uadc0\ud558\uc758 \ubb38\uc11c\uac00 \ud655\uc778\ub418\uc5c8\uc2b5\ub2c8\ub2e4!\"\n\n# 11) Greeting.\nReplace-All $d.Content \"Hi \" \"\uc548\ub155\ud558\uc138\uc694 \"\n\n# 12) Reminder paragraph.\nReplace-All $d.Content \"We\u2019ll be sending out more details about the event soon, including the agenda and travel arrangements, so make sure to check your inbox regularly.\" \"\uc548\uac74\uacfc \uc5ec\ud589 \uc77c\uc815\uc744 \ud3ec\ud568\ud558\uc5ec \ud589\uc0ac\uc5d0 \ub300\ud55c \uc790\uc138\ud55c \ub0b4\uc6a9\uc744 \uace7 \ubcf4\ub0b4\ub4dc\ub9b4 \uc608\uc815\uc774\ub2c8 \ubc1b\uc740 \ud3b8\uc9c0\ud568\uc744 \uc815\uae30\uc801\uc73c\ub85c \ud655\uc778\ud574 \uc8fc\uc2dc\uae30 \ubc14\ub78d\ub2c8\ub2e4.\"\n\n# 13) Country manager contact intro.\nReplace-All $d.Content \"If you have any questions, please contact your country manager, \" \"\uad81\uae08\ud558\uc2e0 \uc0ac\ud56d\uc740, \uadc0\ud558\uc758 \uad6d\uac00 \ub2f4\ub2f9\uc790 \"\n\n# 14) \", at \" -> \"\uc5d0\uac8c \" (between [NAME] and [EMAIL ADDRESS]).\nReplace-All $d.Content \", at \" \"\uc5d0\uac8c \"\n\n# 15) \" or \" between [EMAIL ADDRESS] and [WHATSAPP NO]. The earlier\n# \"live chat or WhatsApp\" occurrence of \" or \" must remain untouched, so\n# scope the search to only the paragraph that contains \"[WHATSAPP NO]\".\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*WHATSAPP NO*\") {\n        Replace-All $p.Range \" or \" \" \ub610\ub294 \"\n    }\n}\n\n# 16) WhatsApp closing sentence.\nReplace-All $d.Content \" (WhatsApp). \" \" (WhatsApp)\uc744 \ud1b5\ud574 \uc5f0\ub77d\ud574 \uc8fc\uc2dc\uae30 \ubc14\ub78d\ub2c8\ub2e4. \"\n\n# 17) The single comment's text.\nfor ($i = 1; $i -le $d.Comments.Count; $i++) {\n    $c = $d.Comments.Item($i)\n    if ($c.Range.Text -eq \"choose either one\") {\n        $c.Range.Text = \"\ub458 \uc911 \ud558\ub098\ub97c \uc120\ud0dd\ud558\uc138\uc694\"\n    }\n}\n"}
